$d = $word.ActiveDocument

$d.Content.Find.Execute("24÷8=3, 0", $true, $false, $false, $false, $false, $true, 1, $false, "84÷2=42, 0", 2) | Out-Null
$d.Content.Find.Execute("94÷8=11, 6", $true, $false, $false, $false, $false, $true, 1, $false, "59÷6=9, 5", 2) | Out-Null
$d.Content.Find.Execute("17÷9=1, 8", $true, $false, $false, $false, $false, $true, 1, $false, "84÷5=16, 4", 2) | Out-Null
$d.Content.Find.Execute("48÷2=24, 0", $true, $false, $false, $false, $false, $true, 1, $false, "31÷3=10, 1", 2) | Out-Null
$d.Content.Find.Execute("78÷7=11, 1", $true, $false, $false, $false, $false, $true, 1, $false, "53÷6=8, 5", 2) | Out-Null
$d.Content.Find.Execute("64÷3=21, 1", $true, $false, $false, $false, $false, $true, 1, $false, "49÷2=24, 1", 2) | Out-Null
$d.Content.Find.Execute("34÷6=5, 4", $true, $false, $false, $false, $false, $true, 1, $false, "49÷2=24, 1", 2) | Out-Null
$d.Content.Find.Execute("63÷9=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "83÷6=13, 5", 2) | Out-Null
$d.Content.Find.Execute("13÷7=1, 6", $true, $false, $false, $false, $false, $true, 1, $false, "64÷7=9, 1", 2) | Out-Null
$d.Content.Find.Execute("38÷7=5, 3", $true, $false, $false, $false, $false, $true, 1, $false, "77÷3=25, 2", 2) | Out-Null
$d.Content.Find.Execute("72÷2=36, 0", $true, $false, $false, $false, $false, $true, 1, $false, "97÷9=10, 7", 2) | Out-Null
$d.Content.Find.Execute("15÷2=7, 1", $true, $false, $false, $false, $false, $true, 1, $false, "14÷8=1, 6", 2) | Out-Null
$d.Content.Find.Execute("21÷9=2, 3", $true, $false, $false, $false, $false, $true, 1, $false, "77÷9=8, 5", 2) | Out-Null
$d.Content.Find.Execute("18÷5=3, 3", $true, $false, $false, $false, $false, $true, 1, $false, "20÷4=5, 0", 2) | Out-Null
$d.Content.Find.Execute("35÷4=8, 3", $true, $false, $false, $false, $false, $true, 1, $false, "86÷8=10, 6", 2) | Out-Null
$d.Content.Find.Execute("78÷9=8, 6", $true, $false, $false, $false, $false, $true, 1, $false, "51÷6=8, 3", 2) | Out-Null
$d.Content.Find.Execute("29÷6=4, 5", $true, $false, $false, $false, $false, $true, 1, $false, "56÷3=18, 2", 2) | Out-Null
$d.Content.Find.Execute("20÷9=2, 2", $true, $false, $false, $false, $false, $true, 1, $false, "92÷4=23, 0", 2) | Out-Null
$d.Content.Find.Execute("39÷3=13, 0", $true, $false, $false, $false, $false, $true, 1, $false, "34÷7=4, 6", 2) | Out-Null
$d.Content.Find.Execute("65÷5=13, 0", $true, $false, $false, $false, $false, $true, 1, $false, "93÷7=13, 2", 2) | Out-Null
$d.Content.Find.Execute("24÷4=6, 0", $true, $false, $false, $false, $false, $true, 1, $false, "60÷3=20, 0", 2) | Out-Null
$d.Content.Find.Execute("37÷2=18, 1", $true, $false, $false, $false, $false, $true, 1, $false, "53÷5=10, 3", 2) | Out-Null
$d.Content.Find.Execute("72÷7=10, 2", $true, $false, $false, $false, $false, $true, 1, $false, "32÷2=16, 0", 2) | Out-Null
$d.Content.Find.Execute("33÷5=6, 3", $true, $false, $false, $false, $false, $true, 1, $false, "75÷3=25, 0", 2) | Out-Null
$d.Content.Find.Execute("21÷8=2, 5", $true, $false, $false, $false, $false, $true, 1, $false, "94÷2=47, 0", 2) | Out-Null
